$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (row => new text)
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "252"
    6  = "0.00078"
    7  = "0.00020"
    8  = "0.00006"
    9  = "0.00030"
    10 = "0.00042"
    11 = "0.00051"
    12 = "0.05179"
}

foreach ($row in $updates.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $updates[$row]
}

# Rows that previously held a tab-separated run of many stats and are
# being collapsed back down to a single summary value.
$t.Cell(44, 1).Range.Text = "99.98"
$t.Cell(45, 1).Range.Text = "0.05"
$t.Cell(46, 1).Range.Text = "230"
